$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2636.077
$ws.Range("I2").Value = 797.1429
$ws.Range("K2").Value = 797.1429
$ws.Range("M2").Value = -684.1429
$ws.Range("H28").Value = 3174.875
$ws.Range("I28").Value = 725
$ws.Range("K28").Value = 725
$ws.Range("M28").Value = -240
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H86").Value = 12829.6
$ws.Range("I86").Value = 7970.857
$ws.Range("K86").Value = 7970.857
$ws.Range("M86").Value = -6847.857
$ws.Range("H89").Value = 12829.6
$ws.Range("I89").Value = 7970.857
$ws.Range("K89").Value = 39854.285
$ws.Range("M89").Value = -34238.285
$ws.Range("H113").Value = 3997.5
$ws.Range("J113").Value = 3997.5
$ws.Range("L113").Value = 3997.5
$ws.Range("N113").Value = -10505.5
$ws.Range("H132").Value = 2392.2144
$ws.Range("I132").Value = 2468.6155
$ws.Range("K132").Value = 7405.8465
$ws.Range("M132").Value = -4875.8465
$ws.Range("H135").Value = 707.2857
$ws.Range("I135").Value = 434.4
$ws.Range("K135").Value = 3909.6
$ws.Range("M135").Value = -1374.6
$ws.Range("H138").Value = 2623.875
$ws.Range("J138").Value = 5000
$ws.Range("L138").Value = 15000
$ws.Range("N138").Value = -25280
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3924.5454
$ws.Range("I32").Value = 4015.9375
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 4015.9375
$ws.Range("L32").Value = 1000
$ws.Range("M32").Value = -3728.9375
$ws.Range("N32").Value = -1574
$ws.Range("H74").Value = 2501.25
$ws.Range("I74").Value = 2501.25
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2501.25
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1627.25
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 2501.25
$ws.Range("I77").Value = 2501.25
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 12506.25
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -8138.25
$ws.Range("N77").ClearContents()
$ws.Range("H102").Value = 3991.25
$ws.Range("I102").Value = 2034
$ws.Range("J102").Value = 7253.3335
$ws.Range("K102").Value = 2034
$ws.Range("L102").Value = 7253.3335
$ws.Range("M102").Value = -412
$ws.Range("N102").Value = -10497.3335
$ws.Range("H122").Value = 1749
$ws.Range("J122").Value = 999
$ws.Range("L122").Value = 2997
$ws.Range("N122").Value = -7897
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 27525.875
$ws.Range("I82").Value = 15051.75
$ws.Range("K82").Value = 15051.75
$ws.Range("M82").Value = -14668.75
$ws.Range("H85").Value = 27525.875
$ws.Range("I85").Value = 15051.75
$ws.Range("K85").Value = 15051.75
$ws.Range("M85").Value = -13725.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 100
$ws.Range("K22").Value = 100
$ws.Range("M22").Value = 250
$ws.Range("H31").Value = 3691.182
$ws.Range("I31").Value = 2810.2
$ws.Range("J31").Value = 5579
$ws.Range("K31").Value = 2810.2
$ws.Range("L31").Value = 5579
$ws.Range("M31").Value = -2515.2
$ws.Range("N31").Value = -6169
$ws.Range("H34").Value = 3691.182
$ws.Range("I34").Value = 2810.2
$ws.Range("J34").Value = 5579
$ws.Range("K34").Value = 2810.2
$ws.Range("L34").Value = 5579
$ws.Range("M34").Value = -2608.2
$ws.Range("N34").Value = -5983
$ws.Range("H62").Value = 4149.8335
$ws.Range("I62").Value = 3225
$ws.Range("K62").Value = 3225
$ws.Range("M62").Value = -2601
$ws.Range("H65").Value = 4149.8335
$ws.Range("I65").Value = 3225
$ws.Range("K65").Value = 16125
$ws.Range("M65").Value = -13005
$ws.Range("H99").Value = 7551.3
$ws.Range("I99").Value = 7612.5557
$ws.Range("K99").Value = 7612.5557
$ws.Range("M99").Value = -6114.5557
$ws.Range("H126").Value = 7551.3
$ws.Range("I126").Value = 7612.5557
$ws.Range("K126").Value = 22837.6671
$ws.Range("M126").Value = -20367.6671
$ws.Range("H132").Value = 2068.0908
$ws.Range("I132").Value = 1340.2858
$ws.Range("K132").Value = 4020.8574
$ws.Range("M132").Value = -1490.8574
$ws.Range("H134").Value = 2975.3333
$ws.Range("I134").Value = 2820.8235
$ws.Range("K134").Value = 8462.4705
$ws.Range("M134").Value = -5927.470499999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 171.5
$ws.Range("J33").Value = 165
$ws.Range("L33").Value = 990
$ws.Range("N33").Value = -1556
$ws.Range("H44").Value = 1333.2727
$ws.Range("I44").Value = 1741.5
$ws.Range("J44").Value = 1100
$ws.Range("K44").Value = 5224.5
$ws.Range("L44").Value = 3300
$ws.Range("M44").Value = -4826.5
$ws.Range("N44").Value = -4096
$ws.Range("H80").Value = 15833.333
$ws.Range("J80").Value = 15833.333
$ws.Range("L80").Value = 47499.999
$ws.Range("N80").Value = -49371.999
$ws.Range("H83").Value = 15833.333
$ws.Range("J83").Value = 15833.333
$ws.Range("L83").Value = 142499.997
$ws.Range("N83").Value = -151859.997
$ws.Range("H92").Value = 487.5
$ws.Range("I92").Value = 483.33334
$ws.Range("K92").Value = 1450.00002
$ws.Range("M92").Value = -202.0000199999999
$ws.Range("H117").Value = 871.5
$ws.Range("I117").Value = 768
$ws.Range("K117").Value = 2304
$ws.Range("M117").Value = 1138
$ws.Range("H129").Value = 603.6667
$ws.Range("I129").Value = 656
$ws.Range("K129").Value = 1968
$ws.Range("M129").Value = 3032
$ws.Range("H139").Value = 10000
$ws.Range("I139").Value = 10000
$ws.Range("J139").Value = 10000
$ws.Range("K139").Value = 30000
$ws.Range("L139").Value = 30000
$ws.Range("M139").Value = -24860
$ws.Range("N139").Value = -40280
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4982.857
$ws.Range("I80").Value = 3300
$ws.Range("K80").Value = 3300
$ws.Range("M80").Value = -2302
$ws.Range("H83").Value = 4982.857
$ws.Range("I83").Value = 3300
$ws.Range("K83").Value = 16500
$ws.Range("M83").Value = -11508
$ws.Range("H102").Value = 526.8182
$ws.Range("I102").Value = 526.8182
$ws.Range("K102").Value = 526.8182
$ws.Range("M102").Value = 1095.1818
$ws.Range("H122").Value = 8336042.5
$ws.Range("I122").Value = 9617317
$ws.Range("K122").Value = 28851951
$ws.Range("M122").Value = -28849501
$ws.Range("H126").Value = 4904.9
$ws.Range("I126").Value = 4894.3335
$ws.Range("K126").Value = 14683.0005
$ws.Range("M126").Value = -12213.0005
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3980.75
$ws.Range("I7").Value = 3980.75
$ws.Range("K7").Value = 3980.75
$ws.Range("M7").Value = -3868.75
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H122").Value = 3860.7144
$ws.Range("I122").Value = 3753.3333
$ws.Range("K122").Value = 11259.9999
$ws.Range("M122").Value = -8809.999899999999
$ws.Range("H126").Value = 3980.75
$ws.Range("I126").Value = 3980.75
$ws.Range("K126").Value = 11942.25
$ws.Range("M126").Value = -9472.25
$ws.Range("H132").Value = 6681.4546
$ws.Range("I132").Value = 5249.3335
$ws.Range("K132").Value = 15748.0005
$ws.Range("M132").Value = -13218.0005
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 23000
$ws.Range("J54").Value = 25750
$ws.Range("L54").Value = 25750
$ws.Range("N54").Value = -26790
$ws.Range("H126").Value = 1867.4
$ws.Range("I126").Value = 1867.4
$ws.Range("K126").Value = 5602.200000000001
$ws.Range("M126").Value = -3132.200000000001
$ws.Range("H136").Value = 2546.261
$ws.Range("I136").Value = 2062.647
$ws.Range("K136").Value = 6187.941
$ws.Range("M136").Value = -3637.941
